$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: append "<sep><label>" to the end of a paragraph's visible text,
# with <label> (e.g. "DONE"/"PENDING") highlighted, while <sep> (e.g. " - ")
# stays un-highlighted. Works around the fact that setting
# Range.HighlightColorIndex on an arbitrary sub-range paints the whole
# paragraph that contains it: we temporarily split the label into its own
# paragraph, highlight that (now single-purpose) paragraph cleanly, then
# delete the paragraph mark to rejoin it with the original paragraph.
# ---------------------------------------------------------------------------
function Add-StatusTag($ParaIndex, $Separator, $Label, $HighlightIndex) {

    $p = $d.Paragraphs($ParaIndex).Range

    # Append the separator (" - ") to the paragraph as plain text.
    $sepPoint = $d.Range($p.End - 1, $p.End - 1)
    $sepPoint.InsertAfter($Separator)

    # Split off a fresh paragraph right after, so the label can be
    # highlighted without touching any neighbouring text.
    $splitPoint = $d.Range($p.End - 1, $p.End - 1)
    $splitPoint.InsertParagraphAfter()

    $labelPara = $d.Paragraphs($ParaIndex + 1).Range
    $labelInsert = $d.Range($labelPara.Start, $labelPara.Start)
    $labelInsert.InsertAfter($Label)

    $labelPara2 = $d.Paragraphs($ParaIndex + 1).Range
    $labelPara2.HighlightColorIndex = $HighlightIndex

    # Delete the paragraph mark that separated the two paragraphs, merging
    # the highlighted label back onto the end of the original paragraph.
    $breakRange = $d.Range($p.End - 1, $p.End)
    $breakRange.Delete()
}

# wdColorIndex constants used below: wdGreen = 4, wdYellow = 7
$wdGreen = 4
$wdYellow = 7

# ---------------------------------------------------------------------------
# Paragraph 1: "...kuziapload kama kawaida,like exiting processes ,but..."
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(" kuziapload kama kawaida,like exiting processes ,but with the following exceptions", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    " kuzi upload kama kawaida, like exiting processes, but with the following exceptions", 2)

# ---------------------------------------------------------------------------
# Paragraph: "2. if not signed,lazima tuwe na ticket # ,email ya TRA showing
# there was an an system failure."
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("2. if not signed,lazima ", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "2. if not signed, lazima ", 2)

$d.Content.Find.Execute("ticket # ,email", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "ticket #, email", 2)

$d.Content.Find.Execute("TRA showing there was an an system failure.", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "TRA showing there was system failure.", 2)

# ---------------------------------------------------------------------------
# Paragraph: "Once all these are available ,then the following steps will
# prevail"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Once all these are available ,then the following steps will prevail", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Once all these are available, then the following steps will prevail:", 2)

# ---------------------------------------------------------------------------
# Item 2: "2. Email with ticket reference will be uploaded,@ezrankayamba to
# add this provision for uploading other docs(attachment)"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("2. Email with ticket reference will be uploaded,@ezrankayamba to add this provision for uploading other docs(attachment)", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "2. Email with ticket reference will be uploaded, @ezrankayamba to add this provision for uploading other docs (attachment)", 2)

# ---------------------------------------------------------------------------
# Item 3: "3 Transactions will be sent to admin for review and approval"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("3 Transactions will be sent to admin for review and approval", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "3. Transactions will be sent to admin for review and approval", 2)

# ---------------------------------------------------------------------------
# Item 4: "4 Once approved,the transaction will be declared complete and
# ready for invoicing "
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("4 Once approved,the transaction will be declared complete and ready for invoicing ", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "4. Once approved, the transaction will be declared complete and ready for invoicing ", 2)

# ---------------------------------------------------------------------------
# Now append the " - DONE" / " - PENDING" status tags to the five numbered
# steps below "... the following steps will prevail:".
# ---------------------------------------------------------------------------

# Re-locate paragraph indices by text, since Find/Replace above did not add
# or remove any paragraphs.
function Get-ParaIndexByStart($Prefix) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $t = $d.Paragraphs($i).Range.Text
        if ($t.StartsWith($Prefix)) {
            return $i
        }
    }
    return -1
}

$idx1 = Get-ParaIndexByStart("1. The DMC will be uploaded as usual")
Add-StatusTag $idx1 " - " "DONE" $wdGreen

$idx2 = Get-ParaIndexByStart("2. Email with ticket reference")
Add-StatusTag $idx2 " - " "PENDING" $wdYellow

$idx3 = Get-ParaIndexByStart("3. Transactions will be sent")
Add-StatusTag $idx3 " - " "DONE" $wdGreen

$idx4 = Get-ParaIndexByStart("4. Once approved")
Add-StatusTag $idx4 "- " "DONE" $wdGreen

$idx5 = Get-ParaIndexByStart("5. All transactions with this nature")
Add-StatusTag $idx5 " - " "PENDING" $wdYellow
